$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "Host Names"
$ws.Range("B2").Value = "cwh-esx09.rackwareinc.lab-cwh-esx09.rackwareinc.lab"
$ws.Range("A2").Value = "Trial 3"
$ws.Range("C2").Value = "hcm-cluster01"
$ws.Range("D2").Value = "cwh-esx09.rackwareinc.lab"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "NA"
$ws.Range("H2").Value = "NA"
$ws.Range("E2").Value = "esx09-datastore3"
$ws.Columns("B:B").AutoFit()
$ws.Range("B2").Value = ""
